# "Test Cases" is the first worksheet in this workbook.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Column D (Runmode) formatting touch-up -----------------------------
# D2:D89 get their style record refreshed (re-apply the thin border that is
# already there) so they line up on the canonical "bordered / no fill"
# style instead of the old redundant one.
$ws.Range("D2:D89").Borders.LineStyle = 1

# --- Row 21: result flips from PASS to SKIP ------------------------------
$ws.Range("E21").Value = "SKIP"

# --- New test cases: B89, B90, B91 --------------------------------------
# Seed formatting (borders / wrap) for the 3 new rows by copying an
# existing fully-formatted row, then overwrite the copied values.
[void]$ws.Range("A89:E89").Copy()
[void]$ws.Range("A90:E90").PasteSpecial(-4122)
$ws.Range("A90").Value = "TestCase_B89"
$ws.Range("B90").Value = "OPQA-575"
$ws.Range("C90").Value = "Verify that DETAILS link is working correctly in record view page of a patent"
$ws.Range("D90").Value = "Y"
$ws.Range("E90").Value = "FAIL"

[void]$ws.Range("A89:E89").Copy()
[void]$ws.Range("A91:E91").PasteSpecial(-4122)
$ws.Range("A91").Value = "TestCase_B90"
$ws.Range("B91").Value = "OPQA-577"
$ws.Range("C91").Value = "Verify that following options get displayed in SORT BY drop down in ARTICLES search results page: a)Relevance b)Times Cited c)Publication Date(Newest) d)Publication Date(Oldest)"
$ws.Range("D91").Value = "Y"
$ws.Range("E91").Value = "PASS"
$ws.Rows.Item(91).RowHeight = 30

[void]$ws.Range("A89:E89").Copy()
[void]$ws.Range("A92:E92").PasteSpecial(-4122)
$ws.Range("A92").Value = "TestCase_B91"
$ws.Range("B92").Value = "OPQA-579"
$ws.Range("C92").Value = "Verify that search results are sorted by RELEVANCE by default in PATENTS search results page"
$ws.Range("D92").Value = "Y"
$ws.Range("E92").Value = "PASS"

# --- Sheet view: scroll so row 88 is near the top, select D2:D92 --------
[void]$ws.Range("D2:D92").Select()
$excel.ActiveWindow.ScrollRow = 88
